$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.138755083084106
$ws.Range("B1").Value = 3.218261241912842
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.558988571166992
$ws.Range("E1").Value = 2.09212589263916
